$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add note about imputation possibility under the "mcig" row (D21)
$ws.Cells.Item(21, 4).Value = "3 of them can be imputed using educage/feg (maybe)"

# Extend the note on the "socstat" row (D22) with additional imputation info
$ws.Cells.Item(22, 4).Value = "NOT the same missings as mcig. 6 of them can be imputed using educage/feg (maybe)"

# Resize columns B, C, D slightly (as a result of re-flowing the new text)
$ws.Columns.Item(2).ColumnWidth = 68.5
$ws.Columns.Item(3).ColumnWidth = 17.75
$ws.Columns.Item(4).ColumnWidth = 84.3

# Update the selected/active cell
$ws.Range("D18").Select()
